$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value while forcing it to remain a literal text
# string even when the text looks like a number or a date (Excel's COM
# layer would otherwise silently convert "300.01" -> 300.01 (number) or
# "2025-09-26" -> a date serial). We briefly mark the cell as Text,
# assign the literal value, then restore the default ("Normal") style
# so no stray formatting is left behind.
function Set-TextValue {
    param($cell, [string]$value)
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value2 = $value
    $rng.Style = "Normal"
}

# --- Row 2 updates -----------------------------------------------------
# expense-1 amount changed from 100 to 0.01 inside the JSON stored in H2
Set-TextValue "H2" '[{"id":"expense-1","type":"transportation","amount":0.01,"banks":["كريدى","مانى فيللوز"]},{"id":"expense-2","type":"transportation","amount":300,"banks":["اسكندرية"]}]'

# totalAmount recalculated accordingly (400 -> 300.01), kept as text
Set-TextValue "I2" "300.01"

# --- New row 4 -----------------------------------------------------------
Set-TextValue "A4" "3f65f14e-3669-4aa4-8858-3a66a39cef52"
$ws.Range("B4").Value2 = 675
Set-TextValue "C4" "كريم خالد محمد محمود"
Set-TextValue "D4" "20أ القاهرة"
Set-TextValue "E4" "2025-09-26"
Set-TextValue "G4" "تتيستسس"
Set-TextValue "H4" '[{"id":"expense-1","type":"transportation","amount":10,"banks":["كريدى","مانى فيللوز"]}]'
Set-TextValue "I4" "10"
Set-TextValue "J4" "2025-09-14T12:02:06.872Z"

# --- Keep the "number stored as text" warning suppressed for the whole
# table now that it spans through row 4 (mirrors the ignoredErrors sqref
# growing from A1:J3 to A1:J4 in the source file).
$errs = $ws.Range("A1:J4").Errors
$errs.Item(3).Ignore = $true

Write-Host "Mission row added and totals updated."
